# ---------------------------------------------------------------------------
# Add a "2022-Q3" sheet (fund-holdings detail) right after "2022-Q2", and
# update the "总计" (totals) sheet with a new row for 2022-Q3, shifting the
# existing history rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q3" worksheet ---------------------------------
# Duplicate the existing "2022-Q2" sheet (so we inherit its sheetPr / page
# setup / column formatting) and place the duplicate right after it, then
# rename + move it in front so the tab order becomes ... 2022-Q2, 2022-Q3 ...
#  -> ... 2022-Q3, 2022-Q2 ...
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"
$q3.Move($q2)

# NOTE: this COM engine resolves a worksheet handle's *position* lazily, so
# after Move() shuffles tab order, previously-grabbed handles like $q2/$q3
# can end up pointing at whatever sheet now sits in their *original* slot.
# Re-fetch the worksheet we actually want to edit by name once the tab
# order has settled.
$ws = $wb.Worksheets.Item("2022-Q3")

# The copied sheet has 13 data rows (rows 2-14); the 2022-Q3 snapshot only
# has 10 data rows (rows 2-11), so drop the trailing 3 rows of old data.
$ws.Range("A12:H14").Delete()

# --- 2. Populate the "2022-Q3" sheet with the snapshot data ----------------

# Columns B (fund code) and D:G (numeric-looking figures kept as text in the
# source data) must stay text so things like leading zeros / trailing zeros
# are preserved instead of being coerced into numbers.
$ws.Range("B2:B11").NumberFormat = "@"
$ws.Range("D2:G11").NumberFormat = "@"

$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "005775"
$ws.Cells.Item(2,3).Value = "中加转型动力灵活配置混合A"
$ws.Cells.Item(2,4).Value = "5.34"
$ws.Cells.Item(2,5).Value = "50.55"
$ws.Cells.Item(2,6).Value = "3.27"
$ws.Cells.Item(2,7).Value = "0.1746"
$ws.Cells.Item(2,8).Value = 3

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "005014"
$ws.Cells.Item(3,3).Value = "泰康景泰回报混合A"
$ws.Cells.Item(3,4).Value = "8.99"
$ws.Cells.Item(3,5).Value = "34.29"
$ws.Cells.Item(3,6).Value = "1.38"
$ws.Cells.Item(3,7).Value = "0.1241"
$ws.Cells.Item(3,8).Value = 8

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "009414"
$ws.Cells.Item(4,3).Value = "中银大健康股票A"
$ws.Cells.Item(4,4).Value = "2.64"
$ws.Cells.Item(4,5).Value = "87.62"
$ws.Cells.Item(4,6).Value = "2.91"
$ws.Cells.Item(4,7).Value = "0.0768"
$ws.Cells.Item(4,8).Value = 10

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "009242"
$ws.Cells.Item(5,3).Value = "中加核心智造混合A"
$ws.Cells.Item(5,4).Value = "1.92"
$ws.Cells.Item(5,5).Value = "61.20"
$ws.Cells.Item(5,6).Value = "4.00"
$ws.Cells.Item(5,7).Value = "0.0768"
$ws.Cells.Item(5,8).Value = 2

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "012072"
$ws.Cells.Item(6,3).Value = "中加喜利回报一年持有期混合C"
$ws.Cells.Item(6,4).Value = "2.21"
$ws.Cells.Item(6,5).Value = "38.64"
$ws.Cells.Item(6,6).Value = "2.37"
$ws.Cells.Item(6,7).Value = "0.0524"
$ws.Cells.Item(6,8).Value = 6

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "012071"
$ws.Cells.Item(7,3).Value = "中加喜利回报一年持有期混合A"
$ws.Cells.Item(7,4).Value = "1.98"
$ws.Cells.Item(7,5).Value = "38.64"
$ws.Cells.Item(7,6).Value = "2.37"
$ws.Cells.Item(7,7).Value = "0.0469"
$ws.Cells.Item(7,8).Value = 6

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "005776"
$ws.Cells.Item(8,3).Value = "中加转型动力灵活配置混合C"
$ws.Cells.Item(8,4).Value = "1.27"
$ws.Cells.Item(8,5).Value = "50.55"
$ws.Cells.Item(8,6).Value = "3.27"
$ws.Cells.Item(8,7).Value = "0.0415"
$ws.Cells.Item(8,8).Value = 3

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "005015"
$ws.Cells.Item(9,3).Value = "泰康景泰回报混合C"
$ws.Cells.Item(9,4).Value = "0.39"
$ws.Cells.Item(9,5).Value = "34.29"
$ws.Cells.Item(9,6).Value = "1.38"
$ws.Cells.Item(9,7).Value = "0.0054"
$ws.Cells.Item(9,8).Value = 8

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "009243"
$ws.Cells.Item(10,3).Value = "中加核心智造混合C"
$ws.Cells.Item(10,4).Value = "0.09"
$ws.Cells.Item(10,5).Value = "61.20"
$ws.Cells.Item(10,6).Value = "4.00"
$ws.Cells.Item(10,7).Value = "0.0036"
$ws.Cells.Item(10,8).Value = 2

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "010321"
$ws.Cells.Item(11,3).Value = "中银大健康股票C"
$ws.Cells.Item(11,4).Value = "0.11"
$ws.Cells.Item(11,5).Value = "87.62"
$ws.Cells.Item(11,6).Value = "2.91"
$ws.Cells.Item(11,7).Value = "0.0032"
$ws.Cells.Item(11,8).Value = 10

# Drop back to the default ("Normal") style so the forced text number format
# doesn't linger as extra explicit cell formatting.
$ws.Range("B2:B11").Style = "Normal"
$ws.Range("D2:G11").Style = "Normal"

# --- 3. Update the "总计" summary sheet -------------------------------------
# Add the 2022-Q3 row at the top of the history and shift every other
# quarter's row down by one (and bump the running index in column A).
# Re-fetch by name (see note above) even though this sheet wasn't moved.
$total = $wb.Worksheets.Item("总计")

# Extend the table by one row first, copying row 8's formatting down to the
# new row 9 so column A keeps its original style.
$total.Range("A8:D8").Copy()
$total.Range("A9:D9").PasteSpecial(-4122)

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 10
$total.Cells.Item(2,4).Value = 0.61

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q2"
$total.Cells.Item(3,3).Value = 13
$total.Cells.Item(3,4).Value = 1.25

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q1"
$total.Cells.Item(4,3).Value = 14
$total.Cells.Item(4,4).Value = 1.39

$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(5,2).Value = "2021-Q4"
$total.Cells.Item(5,3).Value = 15
$total.Cells.Item(5,4).Value = 1.26

$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(6,2).Value = "2021-Q3"
$total.Cells.Item(6,3).Value = 9
$total.Cells.Item(6,4).Value = 0.96

$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(7,2).Value = "2021-Q2"
$total.Cells.Item(7,3).Value = 6
$total.Cells.Item(7,4).Value = 0.58

$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(8,2).Value = "2021-Q1"
$total.Cells.Item(8,3).Value = 18
$total.Cells.Item(8,4).Value = 1.34

$total.Cells.Item(9,1).Value = 7
$total.Cells.Item(9,2).Value = "2020-Q4"
$total.Cells.Item(9,3).Value = 8
$total.Cells.Item(9,4).Value = 0.31

# Leave the "总计" sheet as the active tab, matching the original workbook.
$total.Activate()
